# entregaveis.xlsx — "Atualizacao da Planilha dos Entregaveis e da Qualidade dos entregaveis"
#
# Summary of the edit being reproduced:
#   - A new "Data da conclusao" (completion date) column is inserted right after
#     column B ("Data prevista para conclusao"); the old column C ("Situacao dos
#     entregaveis") slides over to become column D.
#   - All three deliverable rows are marked 100% complete (1 instead of 0/0.9) in
#     the (now) D column, and get an actual completion date in the new C column.
#   - The 2nd deliverable's name changes from "Modulo Documentacao do historico
#     dos pacientes" to "Modulo Relatorio do historico de exames dos pacientes".
#   - The (now orphaned) AutoFilter's hidden "_FilterDatabase" defined name turns
#     into a #REF! error, which happens whenever the filtered range it points to
#     is no longer structurally intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room for the new column: move the old column C (header +
#        3 data cells, formats included) into column D before it gets
#        overwritten. Using copy/paste-formats (rather than Columns.Insert)
#        keeps the sheet's other parts -- e.g. the printer-settings
#        relationship -- untouched, matching how the real edit only ever
#        touches cell content/column width.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)     # xlPasteFormats
$ws.Range("C2:C4").Copy()
$ws.Range("D2:D4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("D1").Value = $ws.Range("C1").Value()
$ws.Range("D2").Value = 1
$ws.Range("D3").Value = 1
$ws.Range("D4").Value = 1

# --- 2. Turn (now free) column C into the new "Data da conclusao" column,
#        copying column B's header/date formatting so the new cells share
#        the existing styles (no new style entries).
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)     # xlPasteFormats
$ws.Range("B2:B4").Copy()
$ws.Range("C2:C4").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C1").Value = "Data da conclusão"
$ws.Range("C2").Value = $ws.Range("B2").Value()   # finished on schedule
$ws.Range("C3").Value = 42114                      # finished a few days late
$ws.Range("C4").Value = $ws.Range("B4").Value()   # finished on schedule

# --- 3. Rename the 2nd deliverable.
$ws.Range("A3").Value = "Modulo Relatório do histórico de exames dos pacientes"

# --- 4. Size the new column the way the workbook ends up sized.
$ws.Columns("D").ColumnWidth = 24.140625

# --- 5. The AutoFilter's hidden defined name no longer resolves to a valid
#        range once the filtered table's shape changed underneath it, so it
#        shows up as a #REF! error, still scoped to the sheet.
$wb.Names.Item(1).RefersTo = "=Entregáveis!#REF!"

Write-Output "entregaveis.xlsx updated"
